# Smart Music Player v1 - update EmotionLinks sheet:
#  - playlist cells now point at .\songs\*.m3u files instead of bare *.mp3 names
#  - a couple of cosmetic view/column-width tweaks from the resave

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: emotion -> playlist file mapping (now .m3u playlists under .\songs\)
$ws.Range("A2").Value = ".\songs\angry.m3u"
$ws.Range("B2").Value = ".\songs\happy.m3u"
$ws.Range("C2").Value = ".\songs\sad.m3u"
$ws.Range("D2").Value = ".\songs\neutral.m3u"

# Column widths tweaked slightly
$ws.Columns.Item(2).ColumnWidth = 23.6666666666667
$ws.Columns.Item(4).ColumnWidth = 19.6666666666667

# Selection moved to D3 (also clears the old frozen top-left scroll position)
$ws.Range("D3").Select() | Out-Null
